$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: label for the "Cantidad" total
$ws.Range("C4").Value = "Total Cantidad"

# Row 5: the total quantity itself, stored as text "3" (same convention the
# sheet already uses for Cantidad/Valor cells). Format the cell as Text
# first so Excel doesn't silently coerce the digit string into a number,
# then strip the formatting again so no stray style sticks to the cell.
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "3"
$ws.Range("C5").ClearFormats()

# Row 6: label for the "Valor" total
$ws.Range("D6").Value = "Total precio"

# Row 7: the total price, again kept as literal text "$70000".
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "$70000"
$ws.Range("D7").ClearFormats()
